# Update the dSF column (F) values for several rows to reflect
# repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    13 = -4
    15 = 0
    18 = -1
    19 = -2
    26 = -2
    28 = -2
    32 = 3
    42 = -1
    57 = 14
    58 = -9
    59 = 8
    60 = -4
    70 = 4
    73 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
